$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.754.81"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.505.94"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.82"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.90"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.29"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("E11").Value = "  +6.15%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "2.899.92"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "2.512.40"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "47.659.04"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.32"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.75"
$ws.Range("E22").Value = "  +7.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.92"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.24"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -4.38%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.142"
$ws.Range("E30").Value = "  +4.01%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.89"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.82"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.02"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.35"
$ws.Range("E41").Value = "  +4.83%  "
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.52"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "1.999.94"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.09"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.17"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.49"
$ws.Range("E51").Value = "  +1.83%  "
